# Append new listing scraped at 2025-12-30 01:25:08 JST into the
# "ランサーズ" sheet, keeping results ordered by priority score (column G).
# The new item (score 45) sorts in between the existing row 7 (score 100)
# and the former rows 8-9 (score 10 each), so it is inserted as the new
# row 8 and the old rows 8-9 shift down to rows 9-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a fresh row at position 8; this shifts the old rows 8 and 9
# (and their formatting, incl. the hyperlink-blue style on column F) down
# to rows 9 and 10.
$ws.Rows.Item(8).Insert()

# Refresh the "取得日時" (fetched-at) timestamp on every existing data row -
# the whole list was re-scraped at the same run time.
$timestamp = "2025-12-30 01:25:08"
$ws.Range("A2").Value = $timestamp
$ws.Range("A3").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A7").Value = $timestamp
$ws.Range("A8").Value = $timestamp
$ws.Range("A9").Value = $timestamp
$ws.Range("A10").Value = $timestamp

# Fill in the newly-inserted row 8 with the new listing's data.
$ws.Range("B8").Value = "【急募】エクセルでの在庫管理システム構築依頼"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5463183"
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = "◇管理"

# Register the hyperlink relationship for the newly-added listing's URL.
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5463183")

# Hyperlinks.Add mints its own (duplicate) cell style; re-apply the sheet's
# existing "Hyperlink" look so F10 keeps the same style index as the other
# link cells instead of growing a redundant one.
$ws.Range("F10").Style = $ws.Range("F9").Style
